# Rename the "partidos" worksheet to "resultados" (tab name only; sheetId / r:id stay the same).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("partidos")
$ws.Name = "resultados"

# Move the active selection on that sheet from C5 to C14.
$ws.Range("C14").Select()
